# Ajout de quelques composants
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: column C label changes from "Composant" to "Composants"
$ws.Range("C2").Value = "Composants"

# Row 8: H8 "?" -> "/"
$ws.Range("H8").Value = "/"

# New row 10: Potentiometre 10K
$ws.Range("C10").Value = "Potentiomètre 10K"
$ws.Range("D10").Value = 0.25
$ws.Range("E10").Value = 2
$ws.Range("F10").Formula = "=D10*E10"
$ws.Range("G10").Value = "Gotronic"
$ws.Range("H10").Formula = '="04601"'
$ws.Range("I10").Value = "https://www.gotronic.fr/art-ajustable-horizontal-10k-8486-117.htm"
$ws.Range("J10").Value = "Fourni par l'école"

# New row 12: Led rouge
$ws.Range("C12").Value = "Led rouge"
$ws.Range("D12").Value = 0.15
$ws.Range("E12").Value = 1
$ws.Range("F12").Formula = "=D12*E12"
$ws.Range("G12").Value = "Gotronic"
$ws.Range("H12").Formula = '="03030"'
$ws.Range("I12").Value = "https://www.gotronic.fr/art-led-5mm-rouge-l51hd-2069.htm"
$ws.Range("J12").Value = "Fourni par l'école"

# Hyperlink on I10
$ws.Hyperlinks.Add($ws.Range("I10"), "https://www.gotronic.fr/art-ajustable-horizontal-10k-8486-117.htm")

# Styles: apply € number format (style index 4 pattern, numFmt 8) to D/F columns for new rows
$ws.Range("D10").Style = $ws.Range("D8").Style
$ws.Range("F10").Style = $ws.Range("F8").Style
$ws.Range("D12").Style = $ws.Range("D8").Style
$ws.Range("F12").Style = $ws.Range("F8").Style
$ws.Range("J10").Style = $ws.Range("J8").Style
$ws.Range("J12").Style = $ws.Range("J8").Style

# Column widths
$ws.Columns("C").ColumnWidth = 16
$ws.Columns("H").ColumnWidth = 16.109375

# View settings
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("G16").Select()
